$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.796.96'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '2.443.05'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '560.85'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.68'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.508'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('E9').Value = '  +7.87%  '
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.332'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('E12').Value = '  -5.01%  '
$ws.Range('E13').Value = '  +4.47%  '
$ws.Range('D14').Value = '68.679.09'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '2.890.69'
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '23.42'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').Value = '2.442.11'
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.60'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '339.27'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.02'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '3.85'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.95'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.36%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.54'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.80'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.84%  '
$ws.Range('D26').Value = '2.567.21'
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.38'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.01'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.29%  '
$ws.Range('D29').Value = '0.0₃0825'
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.20'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.20'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.68%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '433.84'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '159.36'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.00'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.52'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.10%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.38'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.08'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('E44').Value = '  +1.01%  '
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '130.13'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.484'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('E49').Value = '  -1.48%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0924'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('E51').Value = '  +2.44%  '
